# Fix the bug: can't download pic from google
# Remove the two "Runner" rows for Biplav Neupane and Sagar Yonjan from the
# OPERATIONS TEAM sheet, and remove the row for W Dinithi Indrachaya Fernando
# from the STUDENT SUCCESS TEAM sheet.

$wb = $excel.ActiveWorkbook

# OPERATIONS TEAM: delete rows 23 and 24 (Biplav Neupane, Sagar Yonjan)
$wsOps = $wb.Worksheets.Item("OPERATIONS TEAM")
$wsOps.Rows.Item(24).Delete()
$wsOps.Rows.Item(23).Delete()

# STUDENT SUCCESS TEAM: delete row 8 (W Dinithi Indrachaya Fernando)
$wsStudent = $wb.Worksheets.Item("STUDENT SUCCESS TEAM")
$wsStudent.Rows.Item(8).Delete()
